$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Sistema de gerar posters pro evento (PHP to PDF/PNG)" task note (previously
# scheduled in the 44851 week, row 14) is moved up to the 44830 week (row 14 ->
# row 11), replacing "acessibilidade no app". Using Cut/Paste moves both the
# text and the cell's fill/format together, and leaves the source cell (B14)
# blank with its old formatting removed.
$ws.Range("B14").Cut($ws.Range("B11"))

# B14 is left blank; restore its shading to "no fill" explicitly (distinct,
# freshly-built style rather than simply inheriting the old highlighted one).
$ws.Range("B14").Interior.ColorIndex = -4142

# Remove the remaining stale/no-longer-relevant task notes - clearing contents
# only (keeps each cell's existing formatting/fill).
$ws.Range("B12").ClearContents()
$ws.Range("B13").ClearContents()
$ws.Range("B15").ClearContents()

# Update the saved view/selection state.
$ws.Range("D13").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7 | Out-Null
